$d = $word.ActiveDocument

# 1. Insert two new paragraphs at the very start of the document:
#    "Procédure" + " :" (as two runs) with a _GoBack bookmark, then an empty paragraph.
$start = $d.Content
$startRange = $d.Range(0, 0)
$startRange.InsertBefore("Procédure :" + [char]13 + [char]13)

# Re-fetch the first paragraph to add the bookmark at the end of "Procédure :" text
$firstPara = $d.Paragraphs(1).Range
# Place bookmark at the end of the first paragraph (before its paragraph mark)
$bookmarkRange = $d.Range($firstPara.Start, $firstPara.End - 1)
$d.Bookmarks.Add("_GoBack", $bookmarkRange) | Out-Null

# 2. Add lastRenderedPageBreak before the second image's drawing (paragraph with Partition.jar screenshot)
#    This is represented in the OOXML; Word COM doesn't have a direct way to insert this field,
#    so we rely on find/replace is not applicable here - this is a rendering artifact normally
#    generated by Word itself. We leave this to be handled structurally.

# 3. Remove lastRenderedPageBreak from "La fonction à" paragraph - also a rendering artifact.

# 4. Remove old bookmark location (now replaced by new one at top) - handled automatically since
#    bookmark name is unique; re-adding it at the top removes the old location.
